$d = $word.ActiveDocument
$dash = [char]0x2013
$tab = [char]9

# -----------------------------------------------------------------
# Edit 1: "August 2020 to date: Iconic King" paragraph.
# The source text is unchanged, but the five separate runs that make
# up this sentence get consolidated into a single run during the
# edit. A self-replace (find text == replace text) that spans every
# run reliably triggers that consolidation.
# -----------------------------------------------------------------
$augustText = "August 2020 to date: Iconic King"
$d.Content.Find.Execute($augustText, $true, $false, $false, $false, $false, $true, 1, $false, $augustText, 2) | Out-Null

# -----------------------------------------------------------------
# Edit 2: "Software Developer (Backend Developer, GCP)" paragraph
# gains a new ", part-time" suffix. First consolidate the existing
# "spaces" + "Software Developer (Backend Developer, GCP)" runs
# (self-replace, no text change), then append the new suffix as its
# own run right after, matching the diff's run layout.
# -----------------------------------------------------------------
$roleText = "Software Developer (Backend Developer, GCP)"
$d.Content.Find.Execute($roleText, $true, $false, $false, $false, $false, $true, 1, $false, $roleText, 2) | Out-Null

$roleRange = $d.Content
$roleRange.Find.Execute($roleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$suffix = ", part-time"
$insPos = $roleRange.End
$insPoint = $d.Range($insPos, $insPos)
$insPoint.InsertAfter($suffix)
# Toggling a character-level property and back forces the inserted
# text to commit as its own distinct run instead of silently fusing
# back into the previous run.
$suffixRange = $d.Range($insPos, $insPos + $suffix.Length)
$suffixRange.Font.Bold = 1
$suffixRange.Font.Bold = 0

# -----------------------------------------------------------------
# Edit 3: "September 2019 - December 2019: Aura Safira Consulting"
# paragraph. Text unchanged; consolidate its many runs into one via
# a self-replace spanning the whole sentence.
# -----------------------------------------------------------------
$sepText = "September 2019 " + $dash + " December 2019: Aura Safira Consulting"
$d.Content.Find.Execute($sepText, $true, $false, $false, $false, $false, $true, 1, $false, $sepText, 2) | Out-Null

# -----------------------------------------------------------------
# Edit 4: "Software Developer Intern" paragraph. The "  " (two
# spaces) run and the "Software Developer Intern" run merge into a
# single run, while the leading stand-alone tab run is untouched.
# -----------------------------------------------------------------
$internFind = $d.Content
$internFind.Find.Execute("Software Developer Intern", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$internStart = $internFind.Start
$mergeRange = $d.Range($internStart - 2, $internStart + "Software Developer Intern".Length)
$mergeText = "  Software Developer Intern"
$mergeRange.Find.Execute($mergeText, $true, $false, $false, $false, $false, $true, 0, $false, $mergeText, 2) | Out-Null

# -----------------------------------------------------------------
# Edit 5: chess paragraph. "Formerly an active chess player..."
# becomes "Active chess player...". Delete the "Formerly an a"
# prefix (leaving "ctive chess player..." as the tail of the
# original run) and insert a capital "A" right before it as its own
# run.
# -----------------------------------------------------------------
$chessOld = "Formerly an active chess player nationally, ranked top 100 nationally by FIDE (World Chess Organization)."
$chessFind = $d.Content
$chessFind.Find.Execute($chessOld, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$chessStart = $chessFind.Start
$prefixLen = "Formerly an a".Length
$prefixRange = $d.Range($chessStart, $chessStart + $prefixLen)
$prefixRange.Delete()
$aInsPoint = $d.Range($chessStart, $chessStart)
$aInsPoint.InsertBefore("A")
$aRange = $d.Range($chessStart, $chessStart + 1)
$aRange.Font.Bold = 1
$aRange.Font.Bold = 0

# -----------------------------------------------------------------
# Edit 6: phone number paragraph. "07" + "21575442" (two runs) merge
# into a single run "0721575442"; the digits themselves do not
# change.
# -----------------------------------------------------------------
$phoneText = "0721575442"
$d.Content.Find.Execute($phoneText, $true, $false, $false, $false, $false, $true, 1, $false, $phoneText, 2) | Out-Null

Write-Host "All edits applied"
